$wb = $excel.ActiveWorkbook

# --- TableMappings sheet: move table numbers into the table mapping ---
# Column B used to hold a "Table Number" lookup value; it is removed entirely.
# The old column C ("Table Size") shifts left into column B.
# Column A (a running index, previously labelled "Location") is relabelled
# "Table Number" and column B is relabelled "Table Size".
$tm = $wb.Worksheets.Item("TableMappings")
$tm.Columns.Item(2).Delete()
$tm.Range("A1").Value = "Table Number"
$tm.Range("B1").Value = "Table Size"

# --- Companies sheet: the "Table" header is renamed to "Table Number" ---
$co = $wb.Worksheets.Item("Companies")
$co.Range("E1").Value = "Table Number"
$co.Range("E1").Select() | Out-Null

# --- Make TableMappings the active/selected sheet ---
$tm.Activate() | Out-Null
$tm.Range("A1").Select() | Out-Null
